$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description for the "Tank" column (A3) to include an example
$ws.Range("A3").Value = "Name of the Tank. E.g. LP1"

# Move the active cell selection to A4 (matches saved view state)
$ws.Range("A4").Select()
